# The "Recorded By" column (G) on the Session Analysis Results sheet stores
# a comma-separated pair of recorder names/emails. For every data row where
# the primary account "dnasr281@gmail.com" was listed first, swap the two
# entries so "dnasr281@gmail.com" is listed second instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetName = "dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G is 7 ("Recorded By")
    $val = $cell.Value2
    if ($val -ne $null) {
        $parts = $val -split ', ', 2
        if ($parts.Count -eq 2 -and $parts[0] -eq $targetName) {
            $cell.Value2 = $parts[1] + ', ' + $parts[0]
        }
    }
}
